# "Login Seeding - revived"
# The old column C (role labels duplicated into B as "bnb" placeholders)
# is removed; the real role values move into column B, and the selection
# moves to the new B1:B4 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column C entirely - shrinks the used range to A1:B4 and removes
# the old B-column "bnb" placeholder text along with it being overwritten below.
$ws.Columns("C").Delete()

# Column B now carries the actual role values (previously stored in C).
$ws.Range("B1").Value = "captain"
$ws.Range("B2").Value = "captian"
$ws.Range("B3").ClearContents()
$ws.Range("B4").Value = "admin"

# Match the saved selection state (active cell B1, selected range B1:B4).
$ws.Range("B1:B4").Select()
